$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.746.72'
$ws.Range('D3').Value = '3.278.65'
$ws.Range('E3').Value = '  -0.65%  '
$ws.Range('D4').Value = '''0.998'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '''571.00'
$ws.Range('E5').Value = '  -1.31%  '
$ws.Range('D6').Value = '''176.11'
$ws.Range('E6').Value = '  -3.93%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '''0.579'
$ws.Range('E8').Value = '  +1.92%  '
$ws.Range('D9').Value = '3.273.86'
$ws.Range('E9').Value = '  -0.61%  '
$ws.Range('E10').Value = '  -1.65%  '
$ws.Range('D11').Value = '''0.572'
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').Value = '''45.63'
$ws.Range('E12').Value = '  -2.31%  '
$ws.Range('D13').Value = '''0.0000268'
$ws.Range('E13').Value = '  +1.52%  '
$ws.Range('D14').Value = '''696.09'
$ws.Range('E14').Value = '  +9.07%  '
$ws.Range('D15').Value = '3.800.20'
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('E16').Value = '  -1.73%  '
$ws.Range('D17').Value = '66.820.53'
$ws.Range('E17').Value = '  +1.70%  '
$ws.Range('E18').Value = '  +1.03%  '
$ws.Range('D19').Value = '3.280.21'
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('D20').Value = '''17.32'
$ws.Range('E20').Value = '  -2.80%  '
$ws.Range('D21').Value = '''10.72'
$ws.Range('E21').Value = '  -2.36%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('E23').Value = '  -4.17%  '
$ws.Range('D24').Value = '''5.12'
$ws.Range('E24').Value = '  +3.09%  '
$ws.Range('D25').Value = '''98.75'
$ws.Range('E25').Value = '  -2.37%  '
$ws.Range('D26').Value = '''3.88'
$ws.Range('E26').Value = '  -2.06%  '
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('D28').Value = '''9.33'
$ws.Range('E28').Value = '  -0.76%  '
$ws.Range('D29').Value = '''33.09'
$ws.Range('E29').Value = '  +6.89%  '
$ws.Range('D30').Value = '''8.42'
$ws.Range('E30').Value = '  +0.93%  '
$ws.Range('D31').Value = '''6.77'
$ws.Range('E31').Value = '  +4.04%  '
$ws.Range('D32').Value = '''565.85'
$ws.Range('E32').Value = '  -3.76%  '
$ws.Range('D33').Value = '3.896.59'
$ws.Range('E33').Value = '  +1.33%  '
$ws.Range('D34').Value = '''10.83'
$ws.Range('E34').Value = '  -0.50%  '
$ws.Range('E35').Value = '  -1.21%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').Value = '''55.52'
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').Value = '''3.31'
$ws.Range('E38').Value = '  -10.94%  '
$ws.Range('E39').Value = '  +1.56%  '
$ws.Range('E40').Value = '  +0.72%  '
$ws.Range('B41').Value = 'ApeXProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D41').Value = '''3.35'
$ws.Range('E41').Value = '  -1.30%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').Value = '''31.88'
$ws.Range('E42').Value = '  -1.09%  '
$ws.Range('D43').Value = '0.0₃0672'
$ws.Range('E43').Value = '  -3.41%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').Value = '''0.328'
$ws.Range('E44').Value = '  -1.53%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').Value = '''3.00'
$ws.Range('E45').Value = '  -3.82%  '
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('D49').Value = '''2.54'
$ws.Range('E49').Value = '  +1.12%  '
$ws.Range('E50').Value = '  +7.05%  '
$ws.Range('D51').Value = '''129.82'
